# draft-gandhi-spring-stamp-srpm-01.pptx — content edits
#
# Helper: replace the first occurrence of $old with $new inside a
# TextRange, targeting the *exact* substring via .Characters() so the
# underlying run is rewritten in place (same length-of-text edit a user
# would make by selecting the text and retyping it) instead of having
# the whole paragraph collapse into a single new run.
function Replace-InRange {
    param($TextRange, [string]$Old, [string]$New)
    $full = $TextRange.Text
    $idx = $full.IndexOf($Old)
    if ($idx -lt 0) {
        throw "substring not found: $Old"
    }
    $TextRange.Characters($idx + 1, $Old.Length).Text = $New
}

$p = $ppt.ActivePresentation

# Slide 3 ("Requirements and Scope") — rename "SR Policy/Policies" -> "SR Path/Paths"
$s3 = $p.Slides.Item(3)
$shp3 = $s3.Shapes.Item(2)
$tr3 = $shp3.TextFrame.TextRange
Replace-InRange $tr3 "Links and End-to-end P2P/P2MP SR Policies" "Links and End-to-end P2P/P2MP SR Paths"
Replace-InRange $tr3 "Handle ECMP for SR Policies" "Handle ECMP for SR Paths"

# Slide 10 ("Stand-alone LM Message Format for STAMP")
$s10 = $p.Slides.Item(10)

# Drop the trailing period on the last bullet.
$shp10b = $s10.Shapes.Item(3)
$tr10b = $shp10b.TextFrame.TextRange
Replace-InRange $tr10b " is used for LM." " is used for LM"

# ASCII packet diagram: "Session ID" field renamed to "SSID".
$shp10r = $s10.Shapes.Item(5)
$tr10r = $shp10r.TextFrame.TextRange
Replace-InRange $tr10r "|X|B| Reserved  | Block Number  | Session ID                    |" "|X|B| Reserved  | Block Number  | SSID                          |"

# Slide 16 ("Probe Query for SR-MPLS and SRv6 Policy") — title font size 34 -> 32
$s16 = $p.Slides.Item(16)
$shp16 = $s16.Shapes.Item(1)
$shp16.TextFrame.TextRange.Font.Size = 32

# Slide 18 ("ECMP Support for SR Policy" -> "...SR Path")
$s18 = $p.Slides.Item(18)
$shp18t = $s18.Shapes.Item(1)
$tr18t = $shp18t.TextFrame.TextRange
Replace-InRange $tr18t "ECMP Support for SR Policy" "ECMP Support for SR Path"

$shp18b = $s18.Shapes.Item(2)
$tr18b = $shp18b.TextFrame.TextRange
Replace-InRange $tr18b "SR Policy can have ECMP between the ingress and transit nodes, between transit nodes and between transit and egress nodes." "SR Path can have ECMP between the ingress and transit nodes, between transit nodes and between transit and egress nodes."
